$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 28 inherits the same formatting (styles) as row 27 directly above it.
$ws.Range("A27:O27").Copy()
$ws.Range("A28:O28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$row = 28
$ws.Cells.Item($row, 1).Value = 45735.186921296299
$ws.Cells.Item($row, 2).Value = 10
$ws.Cells.Item($row, 3).Value = 6
$ws.Cells.Item($row, 4).Value = 249
$ws.Cells.Item($row, 5).Value = 452
$ws.Cells.Item($row, 6).Value = 428
$ws.Cells.Item($row, 7).Value = 478
$ws.Cells.Item($row, 8).Value = 3437
$ws.Cells.Item($row, 9).Value = 478
$ws.Cells.Item($row, 10).Value = 2026
$ws.Cells.Item($row, 11).Value = 208
$ws.Cells.Item($row, 12).Value = 418
$ws.Cells.Item($row, 13).Value = 30
$ws.Cells.Item($row, 14).Value = 3799
$ws.Cells.Item($row, 15).Value = 4968
